# Apply the diagnostics-log style update:
#  - Update the "old" last row (128)'s saving_timestamp / timestamp_experiment_run
#    values (F128 / G128) to their corrected values.
#  - Append six new log rows (129-134) for further MUTAG / (5)-NN_Classifier_GED
#    functionality-test runs, each with the same experiment metadata columns
#    (A-D) and a pair of timestamp values (F/G). Columns E and H stay blank,
#    consistent with every other row in the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# --- Fix up the existing last row (128) ---
$ws.Cells.Item(128, 6).Value = 45912.65889356482
$ws.Cells.Item(128, 7).Value = 45912.65889335648

# --- Data for the newly appended rows 129-134 ---
$experiment = "Fucntionality_test_MUTAG_with_(5)-NN_Classifier_GED"
$dataset = "MUTAG"
$model = "(5)-NN_Classifier_GED"
$filename = "(5)-NN_Classifier_GED_trained_on_MUTAG.joblib"

$newRows = @(
    @{ Row = 129; F = 45913.53986631944; G = 45913.53986611111 },
    @{ Row = 130; F = 45913.53986631944; G = 45913.53986611111 },
    @{ Row = 131; F = 45913.61776412037; G = 45913.61776390046 },
    @{ Row = 132; F = 45913.61776412037; G = 45913.61776390046 },
    @{ Row = 133; F = 45913.61797832176; G = 45913.61797810185 },
    @{ Row = 134; F = 45913.61797831874; G = 45913.61797810689 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $experiment
    $ws.Cells.Item($r, 2).Value = $dataset
    $ws.Cells.Item($r, 3).Value = $model
    $ws.Cells.Item($r, 4).Value = $filename

    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 6).NumberFormat = $dateFmt

    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 7).NumberFormat = $dateFmt
}
